$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) The hidden "_GoBack" bookmark currently sits inside the third paragraph
#    ("This makes me, as the kids say, Mad AF"), splitting it into two runs.
#    In the target document that bookmark moves to live alone inside a new,
#    otherwise-empty paragraph placed right before the "Buster.lu is a..."
#    paragraph. Remove it from its current spot first so re-adding it later
#    doesn't create a duplicate-named bookmark.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Insert three new blank paragraphs immediately before the first
#    paragraph ("Buster.lu is a...").  They will become, in order:
#      - "Hi guys,"
#      - "This "
#      - an empty paragraph that will hold the relocated bookmark
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1).Range
$firstPara.Collapse(1)
$firstPara.InsertParagraphBefore()
$firstPara.InsertParagraphBefore()
$firstPara.InsertParagraphBefore()

$d.Paragraphs.Item(1).Range.Text = "Hi guys,"
$d.Paragraphs.Item(2).Range.Text = "This "

# ---------------------------------------------------------------------------
# 3) Put the bookmark (and nothing else) into the third of the new
#    paragraphs, matching the target markup exactly:
#      <w:p><w:bookmarkStart .../><w:bookmarkEnd .../></w:p>
# ---------------------------------------------------------------------------
$bookmarkXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs.Item(3).Range.InsertXML($bookmarkXml)

# ---------------------------------------------------------------------------
# 4) Rework the (now un-bookmarked) "This makes me..." paragraph so "This"
#    becomes "That" using the same run split shown in the target markup:
#      "Th" / "at" / " makes me, as the kids say, Mad AF"
# ---------------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i).Range
    $candidateText = $candidate.Text.TrimEnd([char]13, [char]7)
    if ($candidateText -eq "This makes me, as the kids say, Mad AF") {
        $targetPara = $candidate
    }
}

$madAfXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Th</w:t></w:r><w:r><w:t>at</w:t></w:r><w:r><w:t xml:space="preserve"> makes me, as the kids say, Mad AF</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$targetPara.InsertXML($madAfXml)
